$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 398
$newDate = 45968.50834490741

$colA = @('RCA NÃO ATUA', 'RUIDO NA HÉLICE', 'RÁDIO NÃO FUNCIONA/SINTONIZA', 'GRAVAÇÃO FALHA DE PROCESSO', 'LUZ DE JIG', 'FLASH LIGHT FALTANDO COR', 'MAU CONTATO NA LEITURA DO PEN DRIVE', 'BATERIA/PILHA NÃO ATUA', 'APARELHO COM CORPO ESTRANHO', 'FALHA DE INJEÇÃO/SERIGRAFIA', 'SOFTWARE TRAVANDO', 'PONTO BRILHANTE', 'TENSÕES VARIANDO', 'LED/FLASH LIGHT COM LUZ FRACA', 'DISPLAY NÃO ACENDE', 'TENSÃO BAIXA', 'VAZAMENTO DE LUZ', 'RUÍDO NO ÁUDIO DO MIC/FONE/AUX', 'CONTROLE COM POUCA SENSIBILIDADE', 'ABERTURA/GAP', 'NÃO GRAVA/ATUALIZA', 'TERRA ABERTO', 'COATING/SELADOR', 'ALETA NÃO ABRE/FECHA', 'EXCESSO DE DÍGITOS', 'O-CELL COM PELÍCULA', 'LINHA VERTICAL', 'POTÊNCIA MÁXIMA', 'LED NÃO APAGA', 'INTERFERÊNCIA NA IMAGEM', 'SEM VÍDEO NO HDMI', 'LED/DISPLAY PISCANDO', 'SOFTWARE DESATUALIZADO', 'VOLUME MÁXIMO NÃO ATUA', 'LED COM COR DIFERENTE', 'PONTO APAGADO', 'CALÇO/QUADRO APARECENDO', 'APARELHO NÃO LÊ PEN DRIVE', 'FALSA FALHA', 'HI-POT/RIGIDEZ/WI', 'LED COM LUZ INVERTIDA', 'EQUIPAMENTO DE TESTE', 'LINHA HORIZONTAL', 'RUÍDO NO VENTILADOR', 'VENTILADOR NÃO GIRA', 'TECLA DESLOCADA/DANIFICADA', 'MAL MONTADO', 'EMPENADO/AMASSADO', 'TESCON MATERIAL', 'SEM ÁUDIO NO MIC/FONE/AUX', 'LIGA/DESLIGA AUTOMATICAMENTE', 'FALHA VISUAL/MONTAGEM', 'TESCON FALHA DE PROCESSO', 'LÂMPADA FRACA/FORTE', 'CENTELHANDO/RUÍDO', 'SEM VÍDEO NO RF/ANTENA', 'LED NÃO ACENDE', 'FLASH LIGHT NÃO LIGA', 'APARELHO NÃO DESLIGA', 'FUNÇÃO INVERTIDA', 'ÁUDIO OSCILANDO', 'LÂMPADA NÃO ACENDE', 'PLACA EM CURTO', 'CONTROLE NÃO ATUA', 'SEM TENSÃO', 'RISCADO', 'SEM ÁUDIO NO CANAL DO AF', 'SEM ÁUDIO NO ALTO-FALANTE', 'PRATO NÃO GIRA', 'APARELHO NÃO CARREGA', 'SEM IMAGEM/SEM BRILHO', 'MANCHA ESCURA NA TELA', 'CONTAMINAÇÃO', 'TECLA DURA', 'QUEBRADO/DANIFICADO/BATIDO', 'NÃO COMUNICA', 'LÂMPADA NÃO APAGA', 'APARELHO NÃO LIGA', 'VIBRAÇÃO NO ÁUDIO', 'SEM ÁUDIO NO CANAL ESQUERDO/FONE', 'MANCHA', 'RUÍDO NO ÁUDIO', 'VOLUME MÍNIMO NÃO ATUA', 'SEM ÁUDIO NO CANAL DIREITO/FONE', 'SEM SINAL DE WI-FI', 'ESPANADO', 'TECLAS NÃO ATUAM', 'VAZAMENTO DE GÁS', 'VAZAMENTO DE AR', 'FALTANDO', 'FORA DO ESPECIFICADO', 'SEM ÁUDIO GERAL', 'SEM ÁUDIO NO TWEETER', 'BLUETOOTH NÃO FUNCIONA', 'NÃO AQUECE', 'ÁUDIO BAIXO', 'COM REBARBA', 'DESLOCADO', 'FALTANDO DÍGITO NO DISPLAY')
$colB = @('TV', 'TV', 'BBS', 'ARCON', 'BBS', 'CM', 'BBS', 'BBS', 'TV', 'ARCON', 'TW', 'TV', 'BBS', 'BBS', 'CM', 'BBS', 'TW', 'CM', 'TV', 'TV', 'TV', 'MWO', 'ARCON', 'ARCON', 'CM', 'TV', 'TV', 'ARCON', 'TV', 'TV', 'TV', 'BBS', 'TV', 'CM', 'TW', 'TV', 'TV', 'BBS', 'ARCON', 'ARCON', 'CM', 'MWO', 'TV', 'MWO', 'MWO', 'BBS', 'ARCON', 'TV', 'ARCON', 'CM', 'BBS', 'ARCON', 'ARCON', 'MWO', 'MWO', 'TV', 'CM', 'CM', 'MWO', 'MWO', 'BBS', 'MWO', 'TV', 'TV', 'CM', 'TV', 'TW', 'CM', 'MWO', 'BBS', 'TV', 'TV', 'TV', 'TV', 'TV', 'ARCON', 'MWO', 'MWO', 'BBS', 'BBS', 'TW', 'CM', 'CM', 'TV', 'TV', 'ARCON', 'CM', 'ARCON', 'BBS', 'ARCON', 'CM', 'CM', 'BBS', 'BBS', 'MWO', 'BBS', 'TV', 'MWO', 'CM')

for ($i = 0; $i -lt $colA.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $colA[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    $ws.Cells.Item($r, 3).Value2 = $newDate
}

$endRow = $startRow + $colA.Length - 1
$ws.Range("C" + $startRow + ":C" + $endRow).NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output ("Added rows " + $startRow + " to " + $endRow)
